$wb = $excel.ActiveWorkbook

# --- Step 1: insert new "2022-Q1" sheet before "总计" ---
$totalSheet = $wb.Worksheets.Item("总计")
$newSheet = $wb.Worksheets.Add($totalSheet)
$newSheet.Name = "2022-Q1"

# Copy header / index-column styling from an existing per-quarter sheet
$styleSrc = $wb.Worksheets.Item("2021-Q4")
$styleSrc.Range("B1:H1").Copy()
$newSheet.Range("B1:H1").PasteSpecial(-4122)
$styleSrc.Range("A2").Copy()
$newSheet.Range("A2:A13").PasteSpecial(-4122)

# Header row
$newSheet.Range("B1").Value = "基金代码"
$newSheet.Range("C1").Value = "基金名称"
$newSheet.Range("D1").Value = "基金规模"
$newSheet.Range("E1").Value = "股票总仓位"
$newSheet.Range("F1").Value = "仓位占比"
$newSheet.Range("G1").Value = "持有市值(亿元)"
$newSheet.Range("H1").Value = "仓位排名"

# Data rows
# row 2
$newSheet.Range("A2").Value = 0
$newSheet.Range("B2").Value = "'003293"
$newSheet.Range("C2").Value = "易方达科瑞灵活配置混合"
$newSheet.Range("D2").Value = "'34.67"
$newSheet.Range("E2").Value = "'78.17"
$newSheet.Range("F2").Value = "'5.40"
$newSheet.Range("G2").Value = "'1.8722"
$newSheet.Range("H2").Value = 1
# row 3
$newSheet.Range("A3").Value = 1
$newSheet.Range("B3").Value = "'110012"
$newSheet.Range("C3").Value = "易方达科汇灵活配置混合"
$newSheet.Range("D3").Value = "'15.73"
$newSheet.Range("E3").Value = "'75.64"
$newSheet.Range("F3").Value = "'5.44"
$newSheet.Range("G3").Value = "'0.8557"
$newSheet.Range("H3").Value = 1
# row 4
$newSheet.Range("A4").Value = 2
$newSheet.Range("B4").Value = "'010389"
$newSheet.Range("C4").Value = "易方达科益混合A"
$newSheet.Range("D4").Value = "'7.10"
$newSheet.Range("E4").Value = "'92.94"
$newSheet.Range("F4").Value = "'7.22"
$newSheet.Range("G4").Value = "'0.5126"
$newSheet.Range("H4").Value = 3
# row 5
$newSheet.Range("A5").Value = 3
$newSheet.Range("B5").Value = "'011649"
$newSheet.Range("C5").Value = "易方达逆向投资混合A"
$newSheet.Range("D5").Value = "'7.49"
$newSheet.Range("E5").Value = "'85.02"
$newSheet.Range("F5").Value = "'5.66"
$newSheet.Range("G5").Value = "'0.4239"
$newSheet.Range("H5").Value = 2
# row 6
$newSheet.Range("A6").Value = 4
$newSheet.Range("B6").Value = "'002291"
$newSheet.Range("C6").Value = "诺安安鑫灵活配置混合"
$newSheet.Range("D6").Value = "'2.19"
$newSheet.Range("E6").Value = "'81.55"
$newSheet.Range("F6").Value = "'8.50"
$newSheet.Range("G6").Value = "'0.1862"
$newSheet.Range("H6").Value = 2
# row 7
$newSheet.Range("A7").Value = 5
$newSheet.Range("B7").Value = "'011650"
$newSheet.Range("C7").Value = "易方达逆向投资混合C"
$newSheet.Range("D7").Value = "'1.96"
$newSheet.Range("E7").Value = "'85.02"
$newSheet.Range("F7").Value = "'5.66"
$newSheet.Range("G7").Value = "'0.1109"
$newSheet.Range("H7").Value = 2
# row 8
$newSheet.Range("A8").Value = 6
$newSheet.Range("B8").Value = "'010857"
$newSheet.Range("C8").Value = "宝盈祥乐一年持有期混合型证券投资基金A"
$newSheet.Range("D8").Value = "'2.00"
$newSheet.Range("E8").Value = "'36.12"
$newSheet.Range("F8").Value = "'3.07"
$newSheet.Range("G8").Value = "'0.0614"
$newSheet.Range("H8").Value = 4
# row 9
$newSheet.Range("A9").Value = 7
$newSheet.Range("B9").Value = "'008324"
$newSheet.Range("C9").Value = "宝盈祥利稳健配置混合A"
$newSheet.Range("D9").Value = "'1.09"
$newSheet.Range("E9").Value = "'36.53"
$newSheet.Range("F9").Value = "'4.18"
$newSheet.Range("G9").Value = "'0.0456"
$newSheet.Range("H9").Value = 2
# row 10
$newSheet.Range("A10").Value = 8
$newSheet.Range("B10").Value = "'010390"
$newSheet.Range("C10").Value = "易方达科益混合C"
$newSheet.Range("D10").Value = "'0.29"
$newSheet.Range("E10").Value = "'92.94"
$newSheet.Range("F10").Value = "'7.22"
$newSheet.Range("G10").Value = "'0.0209"
$newSheet.Range("H10").Value = 3
# row 11
$newSheet.Range("A11").Value = 9
$newSheet.Range("B11").Value = "'008325"
$newSheet.Range("C11").Value = "宝盈祥利稳健配置混合C"
$newSheet.Range("D11").Value = "'0.43"
$newSheet.Range("E11").Value = "'36.53"
$newSheet.Range("F11").Value = "'4.18"
$newSheet.Range("G11").Value = "'0.0180"
$newSheet.Range("H11").Value = 2
# row 12
$newSheet.Range("A12").Value = 10
$newSheet.Range("B12").Value = "'010858"
$newSheet.Range("C12").Value = "宝盈祥乐一年持有期混合型证券投资基金C"
$newSheet.Range("D12").Value = "'0.14"
$newSheet.Range("E12").Value = "'36.12"
$newSheet.Range("F12").Value = "'3.07"
$newSheet.Range("G12").Value = "'0.0043"
$newSheet.Range("H12").Value = 4
# row 13
$newSheet.Range("A13").Value = 11
$newSheet.Range("B13").Value = "'002952"
$newSheet.Range("C13").Value = "建信多因子量化股票"
$newSheet.Range("D13").Value = "'0.10"
$newSheet.Range("E13").Value = "'91.47"
$newSheet.Range("F13").Value = "'2.76"
$newSheet.Range("G13").Value = "'0.0028"
$newSheet.Range("H13").Value = 8

# --- Step 2: insert a new "2022-Q1" row at top of "总计" summary sheet ---
$ws = $wb.Worksheets.Item("总计")
$ws.Rows("2:2").Insert()
# restore index-column (A) style that Insert() does not carry for the new row
$ws.Cells.Item(3,1).Copy()
$ws.Cells.Item(2,1).PasteSpecial(-4122)
# clear the bold formatting Insert() propagated into the new data cells
$ws.Range("B2:D2").ClearFormats()

$ws.Range("A2").Value = 0
$ws.Range("B2").Value = "2022-Q1"
$ws.Range("C2").Value = 12
$ws.Range("D2").Value = 4.11

# renumber the index column for the rows pushed down
$ws.Range("A3").Value = 1
$ws.Range("A4").Value = 2
$ws.Range("A5").Value = 3

